# Append 8 newly-coded segments (rows 208-215) to Sheet1 -- a re-run of the
# qualitative coding pass ("chen") against the latest set of .mex source
# files, per the commit message "Run through with latest mex files".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$templateRow = 207   # last pre-existing data row
$firstRow    = 208
$lastRow     = 215

# New rows share the coder, code, weight score (0) and coverage-denominator
# pattern of every other row in the sheet; only these per-row fields differ:
$newRows = @(
    @{ Row=208; D="2302"; F="3: 1686"; G="3: 1689"; I="256"; J=4; K=1.2445163498335458E-2; M="1/31/19 13:54:31" }
    @{ Row=209; D="2628"; F="2: 6182"; G="2: 6185"; I="256"; J=4; K=2.46290253063235E-2; M="1/31/19 13:55:01" }
    @{ Row=210; D="2628"; F="2: 6202"; G="2: 6205"; I="256"; J=4; K=2.46290253063235E-2; M="1/31/19 13:55:07" }
    @{ Row=211; D="2628"; F="2: 6217"; G="2: 6220"; I="256"; J=4; K=2.46290253063235E-2; M="1/31/19 13:55:11" }
    @{ Row=212; D="2628"; F="2: 6254"; G="2: 6257"; I="256"; J=4; K=2.46290253063235E-2; M="1/31/19 13:55:15" }
    @{ Row=213; D="2628"; F="2: 6120"; G="2: 6122"; I="32"; J=3; K=1.8471768979742627E-2; M="1/31/19 13:55:31" }
    @{ Row=214; D="2628"; F="2: 6135"; G="2: 6137"; I="32"; J=3; K=1.8471768979742627E-2; M="1/31/19 13:55:35" }
    @{ Row=215; D="3872"; F="2: 3375"; G="2: 3377"; I="‡32"; J=3; K=1.3383894713361589E-2; M="1/31/19 13:56:14" }
)

# 1) Clone the last row's formatting (fills/borders/fonts/number formats)
#    down across the new rows so the table keeps its uniform look.
$ws.Range("A207:M207").Copy() | Out-Null
$ws.Range("A208:M215").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) Fill in the constant columns (color swatch, blank comment/doc-group,
#    code, author) plus the per-row fields captured above.
foreach ($d in $newRows) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value  = "●"   # A: Color swatch
    $ws.Cells.Item($r, 2).Value  = ""        # B: Comment (blank)
    $ws.Cells.Item($r, 3).Value  = ""        # C: Document group (blank)

    # D: Document name -- force Text so the numeric-looking name ("2302",
    #    "2628", "3872", ...) is not coerced into a real number.
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $d.D

    $ws.Cells.Item($r, 5).Value  = "MIC"     # E: Code
    $ws.Cells.Item($r, 6).Value  = $d.F       # F: Begin
    $ws.Cells.Item($r, 7).Value  = $d.G       # G: End
    $ws.Cells.Item($r, 8).Value  = 0          # H: Weight score
    $ws.Cells.Item($r, 9).Value  = $d.I       # I: Segment
    $ws.Cells.Item($r, 10).Value = $d.J       # J: Area
    $ws.Cells.Item($r, 11).Value = $d.K       # K: Coverage %
    $ws.Cells.Item($r, 12).Value = "chen"    # L: Author
    $ws.Cells.Item($r, 13).Value = $d.M       # M: Creation date
}

# 3) NumberFormat="@" above (step 2) swapped column D onto a plain Text
#    style; re-clone D207's real style (fill/border/font, General number
#    format) back over D208:D215 now that the values are safely stored as text.
$ws.Range("D207").Copy() | Out-Null
$ws.Range("D208:D215").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

